$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text format first so Excel stores them as literal strings (matching
# the scraped "price" text cells), not as numeric values.
$textCells = @("D5", "D6", "D9", "D11", "D13", "D15", "D16", "D18", "D21", "D22", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new cell values (both D/price and E/volume columns).
$ws.Range("D2").Value = '26.302.83'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.680.64'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '218.69'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").Value = '0.5272'
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("D9").Value = '0.06442'
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("D11").Value = '0.07516'
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("D12").Value = '1.692.38'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Value = '4.555'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '0.000008515'
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '64.49'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '26.337.06'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '4.941'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '189.88'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '6.216'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").Value = '7.778'
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  +5.77%  '
$ws.Range("D27").Value = '15.87'
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  +11.33%  '
$ws.Range("D29").Value = '1.360'
$ws.Range("E29").Value = '  +5.81%  '
$ws.Range("D30").Value = '1.330'
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("D31").Value = '3.591'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("D32").Value = '3.580'
$ws.Range("E32").Value = '  +1.69%  '
$ws.Range("D33").Value = '1.664'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = '0.6221'
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("D36").Value = '2.396'
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("D37").Value = '2.727'
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("D38").Value = '6.439'
$ws.Range("E38").Value = '  +5.72%  '
$ws.Range("D39").Value = '0.01624'
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").Value = '1.107.65'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '0.8812'
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").Value = '100.76'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '1.834.44'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("E45").Value = '  -3.46%  '
$ws.Range("D46").Value = '56.95'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = '1.012'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = '8.173'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D49").Value = '0.05272'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '0.4301'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  +3.14%  '
